$wb = $excel.ActiveWorkbook

# Sheet "展览" (first sheet) - update column F values for the given rows
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 1091
$ws1.Range("F3").Value  = 400
$ws1.Range("F4").Value  = 1501
$ws1.Range("F5").Value  = 8753
$ws1.Range("F9").Value  = 286
$ws1.Range("F10").Value = 154
$ws1.Range("F12").Value = 3598
$ws1.Range("F13").Value = 50
$ws1.Range("F16").Value = 1421
$ws1.Range("F18").Value = 1124
$ws1.Range("F20").Value = 205
$ws1.Range("F21").Value = 2373
$ws1.Range("F22").Value = 62

# Sheet "全部类型" (fourth sheet) - same update, rows shifted by one for the
# last item because this sheet has an extra row compared to "展览"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 1091
$ws4.Range("F3").Value  = 400
$ws4.Range("F4").Value  = 1501
$ws4.Range("F5").Value  = 8753
$ws4.Range("F9").Value  = 286
$ws4.Range("F10").Value = 154
$ws4.Range("F12").Value = 3598
$ws4.Range("F13").Value = 50
$ws4.Range("F16").Value = 1421
$ws4.Range("F18").Value = 1124
$ws4.Range("F20").Value = 205
$ws4.Range("F21").Value = 2373
$ws4.Range("F23").Value = 62
